# Re-adding requirements 109-136 (Reviews category) that were overwritten.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('B124').Value = '109'
$ws.Range('C124').Value = 'Reviews'
$ws.Range('E124').Value = 'A review can be created for any existing movie through a review creation page.'

$ws.Range('B125').Value = '110'
$ws.Range('C125').Value = 'Reviews'
$ws.Range('D125').Value = '109'
$ws.Range('E125').Value = 'The user must be logged in to create a review.'

$ws.Range('B126').Value = '111'
$ws.Range('C126').Value = 'Reviews'
$ws.Range('D126').Value = '109'
$ws.Range('E126').Value = 'The user can enter a title, movie rating, and review body when creating a review.'

$ws.Range('B127').Value = '112'
$ws.Range('C127').Value = 'Reviews'
$ws.Range('D127').Value = '109'
$ws.Range('E127').Value = 'A new review can be submitted by clicking the submit button on the create review page.'

$ws.Range('B128').Value = '113'
$ws.Range('C128').Value = 'Reviews'
$ws.Range('E128').Value = 'A specific review''s contents can be displayed on a page.'

$ws.Range('B129').Value = '114'
$ws.Range('C129').Value = 'Reviews'
$ws.Range('D129').Value = '113'
$ws.Range('E129').Value = 'The display review page shows the related movie title.'

$ws.Range('B130').Value = '115'
$ws.Range('C130').Value = 'Reviews'
$ws.Range('D130').Value = '113'
$ws.Range('E130').Value = 'The display review page shows the review title, and review body.'

$ws.Range('B131').Value = '116'
$ws.Range('C131').Value = 'Reviews'
$ws.Range('D131').Value = '113'
$ws.Range('E131').Value = 'The display review page shows the review movie rating using star images.'

$ws.Range('B132').Value = '117'
$ws.Range('C132').Value = 'Reviews'
$ws.Range('D132').Value = '113'
$ws.Range('E132').Value = 'The display review page has a user information section.'

$ws.Range('B133').Value = '118'
$ws.Range('C133').Value = 'Reviews'
$ws.Range('D133').Value = '117'
$ws.Range('E133').Value = 'The user information section displays a user''s username, avatar, and when the review was posted.'

$ws.Range('B134').Value = '119'
$ws.Range('C134').Value = 'Reviews'
$ws.Range('D134').Value = '117'
$ws.Range('E134').Value = 'The user information section contains buttons for deleting and editing a review.'

$ws.Range('B135').Value = '120'
$ws.Range('C135').Value = 'Reviews'
$ws.Range('D135').Value = '119'
$ws.Range('E135').Value = 'The delete and edit buttons in the user info section only display for the owning user and comment or review moderators.'

$ws.Range('B136').Value = '121'
$ws.Range('C136').Value = 'Reviews'
$ws.Range('D136').Value = '113'
$ws.Range('E136').Value = 'A comment button is displayed and takes user to a comment creation page when clicked.'

$ws.Range('B137').Value = '122'
$ws.Range('C137').Value = 'Reviews'
$ws.Range('D137').Value = '113'
$ws.Range('E137').Value = 'The comment button is only displayed when the user is logged in.'

$ws.Range('B138').Value = '123'
$ws.Range('C138').Value = 'Reviews'
$ws.Range('D138').Value = '113'
$ws.Range('E138').Value = 'Upvote and downvote buttons are displayed on the review display page.'

$ws.Range('B139').Value = '124'
$ws.Range('C139').Value = 'Reviews'
$ws.Range('D139').Value = '113'
$ws.Range('E139').Value = 'The review''s score is displayed on the review display page.'

$ws.Range('B140').Value = '125'
$ws.Range('C140').Value = 'Reviews'
$ws.Range('D140').Value = '123'
$ws.Range('E140').Value = 'Clicking the upvote or downvote buttons while not logged in displays an error popup.'

$ws.Range('B141').Value = '126'
$ws.Range('C141').Value = 'Reviews'
$ws.Range('D141').Value = '123,124'
$ws.Range('E141').Value = 'Clicking on the upvote or downvote button reflects in the review score.'

$ws.Range('B142').Value = '127'
$ws.Range('C142').Value = 'Reviews'
$ws.Range('D142').Value = '123, 124'
$ws.Range('E142').Value = 'Clicking on the upvote or downvote button sends an ajax request to the server to increment or decrement the review score.'

$ws.Range('B143').Value = '128'
$ws.Range('C143').Value = 'Reviews'
$ws.Range('D143').Value = '123, 124'
$ws.Range('E143').Value = 'Clicking on an already selected upvote or downvote button revokes the user''s vote through an AJAX call.'

$ws.Range('B144').Value = '129'
$ws.Range('C144').Value = 'Reviews'
$ws.Range('D144').Value = '117'
$ws.Range('E144').Value = 'Comments associated with a review are displayed under the review on the review display page.'

$ws.Range('B145').Value = '130'
$ws.Range('C145').Value = 'Reviews'
$ws.Range('D145').Value = '117'
$ws.Range('E145').Value = 'A comment displays its body.'

$ws.Range('B146').Value = '131'
$ws.Range('C146').Value = 'Reviews'
$ws.Range('D146').Value = '130'
$ws.Range('E146').Value = 'A comment has a user section that displays the user''s username, avatar, and comment post date.'

$ws.Range('B147').Value = '132'
$ws.Range('C147').Value = 'Reviews'
$ws.Range('D147').Value = '131'
$ws.Range('E147').Value = 'Edit and delete buttons are displayed in the comment''s user section.'

$ws.Range('B148').Value = '133'
$ws.Range('C148').Value = 'Reviews'
$ws.Range('D148').Value = '132'
$ws.Range('E148').Value = 'The edit and delete buttons for comments only display for the user who created them, and for comment or review moderators.'

$ws.Range('B149').Value = '134'
$ws.Range('C149').Value = 'Reviews'
$ws.Range('D149').Value = '120'
$ws.Range('E149').Value = 'The delete button for a review sends an ajax request to have the review deleted, and then refreshes the page.'

$ws.Range('B150').Value = '135'
$ws.Range('C150').Value = 'Reviews'
$ws.Range('E150').Value = 'If a user attempts to create a review for a non existant movie, a movie not found page is displayed.'

$ws.Range('B151').Value = '136'
$ws.Range('C151').Value = 'Reviews'
$ws.Range('E151').Value = 'If a user attempts to display a review that does no exist, a review not found page is displayed.'

# Update the sheet view to match: scrolled so row 123 is at the top, with C137 selected.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 123
$win.ScrollColumn = 1
[void]$ws.Range('C137').Select()
